# Mapa interactivo PEBCOM - automatic map update
# Inserts a new record (case 5037 / Monroe 3605) as row 32 in the PEBCOM
# sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 32 (pushes old row 32.. down to 33..)
$ws.Rows(32).Insert()

# Columns A, B, D and E hold numeric/date-looking values that must stay
# TEXT (matching the rest of the "Caso" / "F. De Reclamo" / "Comuna" / "OT"
# columns), so force the cell format to Text before writing them.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"

$ws.Range("A32").Value = "5037"
$ws.Range("B32").Value = "3/7/2025"
$ws.Range("C32").Value = "Monroe 3605"
$ws.Range("D32").Value = "12"
$ws.Range("E32").Value = "803825082"
$ws.Range("F32").Value = "PEBCOM"
$ws.Range("G32").Value = "Pendiente"
$ws.Range("H32").Value = "Columna inclinada"
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = "Aplomo"
$ws.Range("K32").Value = "Sin equipos"
$ws.Range("L32").Value = "Pasante"
$ws.Range("M32").Value = -58.471774
$ws.Range("N32").Value = -34.565411
$ws.Range("O32").Value = "Colegiales"
$ws.Range("P32").Value = "Capital Norte"
